$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Robo de equipos...): probability 0.2 -> 0.1 ---
$ws.Range("C5").Value = 0.1

# --- Row 6 (Inundacion...): add "Todas" in Areas Afectadas, update Medidas text ---
$ws.Range("B6").Value = "Todas"
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B6").VerticalAlignment = -4108
$ws.Range("B6").WrapText = $false
$ws.Range("D6").Value = "Ubicar los servidires a 1 m. de altura, a salvo de posibles inundaciones. Ademas equipar la sala con alcantarillas para desagotar rapidamente cualquier fuga de agua."

# --- Remove the "Incendio" risk row entirely (row 7); everything below shifts up ---
$ws.Rows.Item(7).Delete()

# --- New row 7 (was "Corte de energia electrica", now shifted up from row 8) ---
$ws.Range("A7").Value = "Corte de energia electrica debido a fallas por parte del proveedor"
$ws.Range("A7").WrapText = $true
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("A7").VerticalAlignment = -4108

$ws.Range("B7").Value = "Todas"
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108
$ws.Range("B7").WrapText = $false

$ws.Range("C7").Value = 0.3

# --- Update selection to match the authored file ---
$ws.Range("C4").Select()

"done"
